$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 6958169.5
$ws.Cells.Item(43, 9).Value = 26250
$ws.Cells.Item(43, 11).Value = 26250
$ws.Cells.Item(43, 13).Value = -26181
$ws.Cells.Item(64, 8).Value = 4052.4167
$ws.Cells.Item(64, 10).Value = 3804.875
$ws.Cells.Item(64, 12).Value = 3804.875
$ws.Cells.Item(64, 14).Value = -4300.875
$ws.Cells.Item(67, 8).Value = 4052.4167
$ws.Cells.Item(67, 10).Value = 3804.875
$ws.Cells.Item(67, 12).Value = 3804.875
$ws.Cells.Item(67, 14).Value = -5520.875
$ws.Cells.Item(88, 8).Value = 1374369.4
$ws.Cells.Item(88, 9).Value = 972.5
$ws.Cells.Item(88, 10).Value = 1766768.4
$ws.Cells.Item(88, 11).Value = 972.5
$ws.Cells.Item(88, 12).Value = 1766768.4
$ws.Cells.Item(88, 13).Value = -566.5
$ws.Cells.Item(88, 14).Value = -1767580.4
$ws.Cells.Item(91, 8).Value = 1374369.4
$ws.Cells.Item(91, 9).Value = 972.5
$ws.Cells.Item(91, 10).Value = 1766768.4
$ws.Cells.Item(91, 11).Value = 972.5
$ws.Cells.Item(91, 12).Value = 1766768.4
$ws.Cells.Item(91, 13).Value = 431.5
$ws.Cells.Item(91, 14).Value = -1769576.4
$ws.Cells.Item(98, 8).Value = 2603.08
$ws.Cells.Item(98, 9).Value = 2220.2632
$ws.Cells.Item(98, 10).Value = 3815.3333
$ws.Cells.Item(98, 11).Value = 2220.2632
$ws.Cells.Item(98, 12).Value = 3815.3333
$ws.Cells.Item(98, 13).Value = -722.2631999999999
$ws.Cells.Item(98, 14).Value = -6811.3333
$ws.Cells.Item(122, 8).Value = 2603.08
$ws.Cells.Item(122, 9).Value = 2220.2632
$ws.Cells.Item(122, 10).Value = 3815.3333
$ws.Cells.Item(122, 11).Value = 6660.7896
$ws.Cells.Item(122, 12).Value = 11445.9999
$ws.Cells.Item(122, 13).Value = -4210.7896
$ws.Cells.Item(122, 14).Value = -16345.9999
$ws.Cells.Item(125, 8).Value = 1571.2
$ws.Cells.Item(125, 9).Value = 1634.6666
$ws.Cells.Item(125, 10).Value = 1000
$ws.Cells.Item(125, 11).Value = 14711.9994
$ws.Cells.Item(125, 12).Value = 9000
$ws.Cells.Item(125, 13).Value = -12251.9994
$ws.Cells.Item(125, 14).Value = -13920
$ws.Cells.Item(131, 8).Value = 1770
$ws.Cells.Item(131, 9).Value = 1770
$ws.Cells.Item(131, 11).Value = 5310
$ws.Cells.Item(131, 13).Value = -270
$ws.Cells.Item(132, 8).Value = 9011324
$ws.Cells.Item(132, 9).Value = 9261467
$ws.Cells.Item(132, 11).Value = 27784401
$ws.Cells.Item(132, 13).Value = -27781871
$ws.Cells.Item(137, 8).Value = 1386.4166
$ws.Cells.Item(137, 9).Value = 1304.0625
$ws.Cells.Item(137, 10).Value = 1551.125
$ws.Cells.Item(137, 11).Value = 3912.1875
$ws.Cells.Item(137, 12).Value = 4653.375
$ws.Cells.Item(137, 13).Value = -1362.1875
$ws.Cells.Item(137, 14).Value = -9753.375
$ws.Cells.Item(141, 8).Value = 895
$ws.Cells.Item(141, 9).Value = 895
$ws.Cells.Item(141, 11).Value = 2685
$ws.Cells.Item(141, 13).Value = 2495

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1449
$ws.Cells.Item(45, 9).Value = 1662.3
$ws.Cells.Item(45, 11).Value = 1662.3
$ws.Cells.Item(45, 13).Value = -1285.3
$ws.Cells.Item(61, 8).Value = 1029.4546
$ws.Cells.Item(61, 9).Value = 832.45
$ws.Cells.Item(61, 11).Value = 832.45
$ws.Cells.Item(61, 13).Value = -620.45
$ws.Cells.Item(74, 8).Value = 1282.8422
$ws.Cells.Item(74, 9).Value = 958.53845
$ws.Cells.Item(74, 11).Value = 958.53845
$ws.Cells.Item(74, 13).Value = -84.53845000000001
$ws.Cells.Item(77, 8).Value = 1282.8422
$ws.Cells.Item(77, 9).Value = 958.53845
$ws.Cells.Item(77, 11).Value = 4792.69225
$ws.Cells.Item(77, 13).Value = -424.6922500000001
$ws.Cells.Item(122, 8).Value = 2796
$ws.Cells.Item(122, 9).Value = 2995
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 8985
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -6535
$ws.Cells.Item(122, 14).Value = -10900
$ws.Cells.Item(132, 8).Value = 3008
$ws.Cells.Item(132, 9).Value = 3619.4285
$ws.Cells.Item(132, 11).Value = 10858.2855
$ws.Cells.Item(132, 13).Value = -8328.2855
$ws.Cells.Item(136, 8).Value = 1029.4546
$ws.Cells.Item(136, 9).Value = 832.45
$ws.Cells.Item(136, 11).Value = 2497.35
$ws.Cells.Item(136, 13).Value = 52.64999999999964

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3040.8918
$ws.Cells.Item(86, 9).Value = 3298.5186
$ws.Cells.Item(86, 11).Value = 3298.5186
$ws.Cells.Item(86, 13).Value = -2175.5186
$ws.Cells.Item(89, 8).Value = 3040.8918
$ws.Cells.Item(89, 9).Value = 3298.5186
$ws.Cells.Item(89, 11).Value = 16492.593
$ws.Cells.Item(89, 13).Value = -10876.593
$ws.Cells.Item(105, 8).Value = 58825150
$ws.Cells.Item(105, 9).Value = 83334776
$ws.Cells.Item(105, 10).Value = 2039.4
$ws.Cells.Item(105, 11).Value = 83334776
$ws.Cells.Item(105, 12).Value = 2039.4
$ws.Cells.Item(105, 13).Value = -83333029
$ws.Cells.Item(105, 14).Value = -5533.4
$ws.Cells.Item(134, 8).Value = 7363.05
$ws.Cells.Item(134, 9).Value = 947.3570999999999
$ws.Cells.Item(134, 10).Value = 22333
$ws.Cells.Item(134, 11).Value = 2842.0713
$ws.Cells.Item(134, 12).Value = 66999
$ws.Cells.Item(134, 13).Value = -307.0712999999996
$ws.Cells.Item(134, 14).Value = -72069

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1914
$ws.Cells.Item(31, 9).Value = 1849.5
$ws.Cells.Item(31, 11).Value = 1849.5
$ws.Cells.Item(31, 13).Value = -1554.5
$ws.Cells.Item(34, 8).Value = 1914
$ws.Cells.Item(34, 9).Value = 1849.5
$ws.Cells.Item(34, 11).Value = 1849.5
$ws.Cells.Item(34, 13).Value = -1647.5
$ws.Cells.Item(58, 8).Value = 948.0278
$ws.Cells.Item(58, 9).Value = 724.25
$ws.Cells.Item(58, 10).Value = 1731.25
$ws.Cells.Item(58, 11).Value = 724.25
$ws.Cells.Item(58, 12).Value = 1731.25
$ws.Cells.Item(58, 13).Value = -521.25
$ws.Cells.Item(58, 14).Value = -2137.25
$ws.Cells.Item(62, 8).Value = 10003915
$ws.Cells.Item(62, 9).Value = 4164.353
$ws.Cells.Item(62, 10).Value = 66669170
$ws.Cells.Item(62, 11).Value = 4164.353
$ws.Cells.Item(62, 12).Value = 66669170
$ws.Cells.Item(62, 13).Value = -3540.353
$ws.Cells.Item(62, 14).Value = -66670418
$ws.Cells.Item(65, 8).Value = 10003915
$ws.Cells.Item(65, 9).Value = 4164.353
$ws.Cells.Item(65, 10).Value = 66669170
$ws.Cells.Item(65, 11).Value = 20821.765
$ws.Cells.Item(65, 12).Value = 333345850
$ws.Cells.Item(65, 13).Value = -17701.765
$ws.Cells.Item(65, 14).Value = -333352090
$ws.Cells.Item(132, 8).Value = 6100.7407
$ws.Cells.Item(132, 9).Value = 8244.134
$ws.Cells.Item(132, 11).Value = 24732.402
$ws.Cells.Item(132, 13).Value = -22202.402
$ws.Cells.Item(136, 8).Value = 948.0278
$ws.Cells.Item(136, 9).Value = 724.25
$ws.Cells.Item(136, 10).Value = 1731.25
$ws.Cells.Item(136, 11).Value = 2172.75
$ws.Cells.Item(136, 12).Value = 5193.75
$ws.Cells.Item(136, 13).Value = 377.25
$ws.Cells.Item(136, 14).Value = -10293.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 90971.82000000001
$ws.Cells.Item(8, 9).Value = 90971.82000000001
$ws.Cells.Item(8, 11).Value = 272915.46
$ws.Cells.Item(8, 13).Value = -272776.46
$ws.Cells.Item(109, 8).Value = 58269.5
$ws.Cells.Item(109, 9).Value = 143650.14
$ws.Cells.Item(109, 10).Value = 3936.3635
$ws.Cells.Item(109, 11).Value = 430950.42
$ws.Cells.Item(109, 12).Value = 11809.0905
$ws.Cells.Item(109, 13).Value = -429910.42
$ws.Cells.Item(109, 14).Value = -13889.0905
$ws.Cells.Item(131, 8).Value = 1182.8485
$ws.Cells.Item(131, 10).Value = 1232.0646
$ws.Cells.Item(131, 12).Value = 3696.1938
$ws.Cells.Item(131, 14).Value = -13776.1938

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 4444
$ws.Cells.Item(39, 10).Value = 4444
$ws.Cells.Item(39, 12).Value = 4444
$ws.Cells.Item(39, 14).Value = -5508
$ws.Cells.Item(102, 8).Value = 2303.75
$ws.Cells.Item(102, 9).Value = 2390.6667
$ws.Cells.Item(102, 11).Value = 2390.6667
$ws.Cells.Item(102, 13).Value = -768.6667000000002
$ws.Cells.Item(107, 8).Value = 524.6087
$ws.Cells.Item(107, 9).Value = 749.7
$ws.Cells.Item(107, 10).Value = 351.46155
$ws.Cells.Item(107, 11).Value = 749.7
$ws.Cells.Item(107, 12).Value = 351.46155
$ws.Cells.Item(107, 13).Value = 1170.3
$ws.Cells.Item(107, 14).Value = -4191.46155
$ws.Cells.Item(122, 8).Value = 252236
$ws.Cells.Item(122, 9).Value = 4850
$ws.Cells.Item(122, 10).Value = 375929
$ws.Cells.Item(122, 11).Value = 14550
$ws.Cells.Item(122, 12).Value = 1127787
$ws.Cells.Item(122, 13).Value = -12100
$ws.Cells.Item(122, 14).Value = -1132687
$ws.Cells.Item(132, 8).Value = 2632.6667
$ws.Cells.Item(132, 9).Value = 2283.7273
$ws.Cells.Item(132, 11).Value = 6851.1819
$ws.Cells.Item(132, 13).Value = -4321.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3125
$ws.Cells.Item(7, 9).Value = 3000
$ws.Cells.Item(7, 10).Value = 3500
$ws.Cells.Item(7, 11).Value = 3000
$ws.Cells.Item(7, 12).Value = 3500
$ws.Cells.Item(7, 13).Value = -2888
$ws.Cells.Item(7, 14).Value = -3724
$ws.Cells.Item(46, 8).Value = 8969.200000000001
$ws.Cells.Item(46, 9).Value = 900
$ws.Cells.Item(46, 11).Value = 900
$ws.Cells.Item(46, 13).Value = -712
$ws.Cells.Item(122, 8).Value = 31253462
$ws.Cells.Item(122, 9).Value = 62503720
$ws.Cells.Item(122, 10).Value = 3202.5
$ws.Cells.Item(122, 11).Value = 187511160
$ws.Cells.Item(122, 12).Value = 9607.5
$ws.Cells.Item(122, 13).Value = -187508710
$ws.Cells.Item(122, 14).Value = -14507.5
$ws.Cells.Item(126, 8).Value = 3125
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 10).Value = 3500
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 12).Value = 10500
$ws.Cells.Item(126, 13).Value = -6530
$ws.Cells.Item(126, 14).Value = -15440
$ws.Cells.Item(130, 8).Value = 74995
$ws.Cells.Item(130, 10).Value = 74995
$ws.Cells.Item(130, 12).Value = 74995
$ws.Cells.Item(130, 14).Value = -85035
$ws.Cells.Item(132, 8).Value = 20052.852
$ws.Cells.Item(132, 9).Value = 1288.4193
$ws.Cells.Item(132, 10).Value = 45344.043
$ws.Cells.Item(132, 11).Value = 3865.2579
$ws.Cells.Item(132, 12).Value = 136032.129
$ws.Cells.Item(132, 13).Value = -1335.2579
$ws.Cells.Item(132, 14).Value = -141092.129

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 10002216
$ws.Cells.Item(122, 9).Value = 12383481
$ws.Cells.Item(122, 10).Value = 903.6
$ws.Cells.Item(122, 11).Value = 37150443
$ws.Cells.Item(122, 12).Value = 2710.8
$ws.Cells.Item(122, 13).Value = -37147993
$ws.Cells.Item(122, 14).Value = -7610.8
$ws.Cells.Item(126, 8).Value = 100001010
$ws.Cells.Item(126, 9).Value = 100001010
$ws.Cells.Item(126, 11).Value = 300003030
$ws.Cells.Item(126, 13).Value = -300000560
$ws.Cells.Item(132, 8).Value = 2951.08
$ws.Cells.Item(132, 9).Value = 2437.7778
$ws.Cells.Item(132, 11).Value = 7313.3334
$ws.Cells.Item(132, 13).Value = -4783.3334
